$d = $word.ActiveDocument
$d.Content.Find.Execute("aGoo", $true, $false, $false, $false, $false, $true, 1, $false, "Goo", 2)
